# Excel single data read SoftAssert
$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Logindata")

# Add a new worksheet named "Item" right after the existing "Logindata" sheet
$wsItem = $wb.Worksheets.Add([Type]::Missing, $wsLogin)
$wsItem.Name = "Item"

# Populate the new sheet with the item-rate data
$wsItem.Range("A1").Value = "ItemRate"
$wsItem.Range("A2").Value = 12

# Match the author's recorded selection / active-tab state
[void]$wsItem.Range("D4").Select()
